# remove institution from SyntheticDataPipeline
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Remove the "Institution" sheet entirely.
# ---------------------------------------------------------------------
$wsInstitution = $wb.Worksheets.Item("Institution")
$wsInstitution.Delete()

# ---------------------------------------------------------------------
# 2. License sheet: append " License" to the BSD / MIT titles.
# ---------------------------------------------------------------------
$wsLicense = $wb.Worksheets.Item("License")
$wsLicense.Range("G2").Value = "BSD License"
$wsLicense.Range("G38").Value = "MIT License"

# ---------------------------------------------------------------------
# 3. RightsStatement sheet: rotate the "note" column text between rows.
# ---------------------------------------------------------------------
$wsRights = $wb.Worksheets.Item("RightsStatement")
$wsRights.Range("E3").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."
$wsRights.Range("E4").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$wsRights.Range("E6").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."
$wsRights.Range("E7").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$wsRights.Range("E8").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."
$wsRights.Range("E9").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
$wsRights.Range("E11").Value = "You may need to obtain other permissions for your intended use. For example, other rights such as publicity, privacy or moral rights may limit how you may use the material."

# ---------------------------------------------------------------------
# 4. Image sheet: re-label Institution-based rows as Collection / Freestanding
#    work rows, then drop the now-unused trailing rows (258-275).
# ---------------------------------------------------------------------
$wsImage = $wb.Worksheets.Item("Image")

$wsImage.Range("A232").Value = "https://place-hold.it/1000x1000?text=Collection0Work0Image0"
$wsImage.Range("C232").Value = "Collection0Work0 image 0"
$wsImage.Range("D232").Value = "http://example.com/collection0/work0"
$wsImage.Range("H232").Value = "Collection0Work0 image 0 rights holder"

$wsImage.Range("A233").Value = "https://place-hold.it/1000x1000?text=Collection0Work0Image1"
$wsImage.Range("C233").Value = "Collection0Work0 image 1"
$wsImage.Range("D233").Value = "http://example.com/collection0/work0"
$wsImage.Range("H233").Value = "Collection0Work0 image 1 rights holder"

$wsImage.Range("A234").Value = "https://place-hold.it/1000x1000?text=Collection0Work1Image0"
$wsImage.Range("C234").Value = "Collection0Work1 image 0"
$wsImage.Range("D234").Value = "http://example.com/collection0/work1"
$wsImage.Range("H234").Value = "Collection0Work1 image 0 rights holder"

$wsImage.Range("A235").Value = "https://place-hold.it/1000x1000?text=Collection0Work1Image1"
$wsImage.Range("C235").Value = "Collection0Work1 image 1"
$wsImage.Range("D235").Value = "http://example.com/collection0/work1"
$wsImage.Range("H235").Value = "Collection0Work1 image 1 rights holder"

$wsImage.Range("A236").Value = "https://place-hold.it/1000x1000?text=Collection0Work2Image0"
$wsImage.Range("C236").Value = "Collection0Work2 image 0"
$wsImage.Range("D236").Value = "http://example.com/collection0/work2"
$wsImage.Range("H236").Value = "Collection0Work2 image 0 rights holder"

$wsImage.Range("A237").Value = "https://place-hold.it/1000x1000?text=Collection0Work2Image1"
$wsImage.Range("C237").Value = "Collection0Work2 image 1"
$wsImage.Range("D237").Value = "http://example.com/collection0/work2"
$wsImage.Range("H237").Value = "Collection0Work2 image 1 rights holder"

$wsImage.Range("A238").Value = "https://place-hold.it/1000x1000?text=Collection0Work3Image0"
$wsImage.Range("C238").Value = "Collection0Work3 image 0"
$wsImage.Range("D238").Value = "http://example.com/collection0/work3"
$wsImage.Range("H238").Value = "Collection0Work3 image 0 rights holder"

$wsImage.Range("A239").Value = "https://place-hold.it/1000x1000?text=Collection0Work3Image1"
$wsImage.Range("C239").Value = "Collection0Work3 image 1"
$wsImage.Range("D239").Value = "http://example.com/collection0/work3"
$wsImage.Range("H239").Value = "Collection0Work3 image 1 rights holder"

$wsImage.Range("A240").Value = "https://place-hold.it/1000x1000?text=Collection1Image0"
$wsImage.Range("C240").Value = "Collection1 image 0"
$wsImage.Range("D240").Value = "http://example.com/collection1"
$wsImage.Range("H240").Value = "Collection1 image 0 rights holder"

$wsImage.Range("A241").Value = "https://place-hold.it/1000x1000?text=Collection1Image1"
$wsImage.Range("C241").Value = "Collection1 image 1"
$wsImage.Range("D241").Value = "http://example.com/collection1"
$wsImage.Range("H241").Value = "Collection1 image 1 rights holder"

$wsImage.Range("A242").Value = "https://place-hold.it/1000x1000?text=Collection1Work4Image0"
$wsImage.Range("C242").Value = "Collection1Work4 image 0"
$wsImage.Range("D242").Value = "http://example.com/collection1/work4"
$wsImage.Range("H242").Value = "Collection1Work4 image 0 rights holder"

$wsImage.Range("A243").Value = "https://place-hold.it/1000x1000?text=Collection1Work4Image1"
$wsImage.Range("C243").Value = "Collection1Work4 image 1"
$wsImage.Range("D243").Value = "http://example.com/collection1/work4"
$wsImage.Range("H243").Value = "Collection1Work4 image 1 rights holder"

$wsImage.Range("A244").Value = "https://place-hold.it/1000x1000?text=Collection1Work5Image0"
$wsImage.Range("C244").Value = "Collection1Work5 image 0"
$wsImage.Range("D244").Value = "http://example.com/collection1/work5"
$wsImage.Range("H244").Value = "Collection1Work5 image 0 rights holder"

$wsImage.Range("A245").Value = "https://place-hold.it/1000x1000?text=Collection1Work5Image1"
$wsImage.Range("C245").Value = "Collection1Work5 image 1"
$wsImage.Range("D245").Value = "http://example.com/collection1/work5"
$wsImage.Range("H245").Value = "Collection1Work5 image 1 rights holder"

$wsImage.Range("A246").Value = "https://place-hold.it/1000x1000?text=Collection1Work6Image0"
$wsImage.Range("C246").Value = "Collection1Work6 image 0"
$wsImage.Range("D246").Value = "http://example.com/collection1/work6"
$wsImage.Range("H246").Value = "Collection1Work6 image 0 rights holder"

$wsImage.Range("A247").Value = "https://place-hold.it/1000x1000?text=Collection1Work6Image1"
$wsImage.Range("C247").Value = "Collection1Work6 image 1"
$wsImage.Range("D247").Value = "http://example.com/collection1/work6"
$wsImage.Range("H247").Value = "Collection1Work6 image 1 rights holder"

$wsImage.Range("A248").Value = "https://place-hold.it/1000x1000?text=Collection1Work7Image0"
$wsImage.Range("C248").Value = "Collection1Work7 image 0"
$wsImage.Range("D248").Value = "http://example.com/collection1/work7"
$wsImage.Range("H248").Value = "Collection1Work7 image 0 rights holder"

$wsImage.Range("A249").Value = "https://place-hold.it/1000x1000?text=Collection1Work7Image1"
$wsImage.Range("C249").Value = "Collection1Work7 image 1"
$wsImage.Range("D249").Value = "http://example.com/collection1/work7"
$wsImage.Range("H249").Value = "Collection1Work7 image 1 rights holder"

$wsImage.Range("A250").Value = "https://place-hold.it/1000x1000?text=FreestandingWork8Image0"
$wsImage.Range("C250").Value = "FreestandingWork8 image 0"
$wsImage.Range("D250").Value = "http://example.com/freestandingwork8"
$wsImage.Range("H250").Value = "FreestandingWork8 image 0 rights holder"

$wsImage.Range("A251").Value = "https://place-hold.it/1000x1000?text=FreestandingWork8Image1"
$wsImage.Range("C251").Value = "FreestandingWork8 image 1"
$wsImage.Range("D251").Value = "http://example.com/freestandingwork8"
$wsImage.Range("H251").Value = "FreestandingWork8 image 1 rights holder"

$wsImage.Range("A252").Value = "https://place-hold.it/1000x1000?text=FreestandingWork9Image0"
$wsImage.Range("C252").Value = "FreestandingWork9 image 0"
$wsImage.Range("D252").Value = "http://example.com/freestandingwork9"
$wsImage.Range("H252").Value = "FreestandingWork9 image 0 rights holder"

$wsImage.Range("A253").Value = "https://place-hold.it/1000x1000?text=FreestandingWork9Image1"
$wsImage.Range("C253").Value = "FreestandingWork9 image 1"
$wsImage.Range("D253").Value = "http://example.com/freestandingwork9"
$wsImage.Range("H253").Value = "FreestandingWork9 image 1 rights holder"

$wsImage.Range("A254").Value = "https://place-hold.it/1000x1000?text=FreestandingWork10Image0"
$wsImage.Range("C254").Value = "FreestandingWork10 image 0"
$wsImage.Range("D254").Value = "http://example.com/freestandingwork10"
$wsImage.Range("H254").Value = "FreestandingWork10 image 0 rights holder"

$wsImage.Range("A255").Value = "https://place-hold.it/1000x1000?text=FreestandingWork10Image1"
$wsImage.Range("C255").Value = "FreestandingWork10 image 1"
$wsImage.Range("D255").Value = "http://example.com/freestandingwork10"
$wsImage.Range("H255").Value = "FreestandingWork10 image 1 rights holder"

$wsImage.Range("A256").Value = "https://place-hold.it/1000x1000?text=FreestandingWork11Image0"
$wsImage.Range("C256").Value = "FreestandingWork11 image 0"
$wsImage.Range("D256").Value = "http://example.com/freestandingwork11"
$wsImage.Range("H256").Value = "FreestandingWork11 image 0 rights holder"

$wsImage.Range("A257").Value = "https://place-hold.it/1000x1000?text=FreestandingWork11Image1"
$wsImage.Range("C257").Value = "FreestandingWork11 image 1"
$wsImage.Range("D257").Value = "http://example.com/freestandingwork11"
$wsImage.Range("H257").Value = "FreestandingWork11 image 1 rights holder"

# Drop the trailing rows that used to describe Institution1's extra
# collection/shared/freestanding works (now out of range).
$wsImage.Range("A258:I275").EntireRow.Delete()

# ---------------------------------------------------------------------
# 5. Person sheet: swap the wikidata / wikipedia relation links.
# ---------------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("F2").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
$wsPerson.Range("F3").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
$wsPerson.Range("F4").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
$wsPerson.Range("F5").Value = "http://www.wikidata.org/entity/Q7251"
